$wb = $excel.ActiveWorkbook

$wsSummary   = $wb.Worksheets.Item("Summary")
$wsRepay     = $wb.Worksheets.Item("Repayment schedule")
$wsTrans     = $wb.Worksheets.Item("Transactions")

# --- Repayment schedule sheet: K2 / L2 change from 0 to 100 ---
$wsRepay.Range("K2").Value = 100
$wsRepay.Range("L2").Value = 100

# --- Transactions sheet: drop the "Accrual" rows, keep the rest, renumber IDs ---
# Before (rows 2-7):
#   2: 116  Accrual
#   3: 115  Accrual
#   4: 113  Disbursement
#   5: 114  Accrual
#   6: 112  Repayment (at time of disbursement)
#   7: 111  Disbursement
# After deleting the Accrual rows (2,3,5) the remaining rows 4,6,7 shift up to 2,3,4.
$wsTrans.Rows("5:5").Delete()
$wsTrans.Rows("2:3").Delete()

# Renumber the transaction IDs on the now-shifted rows
$wsTrans.Range("A2").Value = 6555
$wsTrans.Range("A3").Value = 6553
$wsTrans.Range("A4").Value = 6552

# --- Update remembered cell selections on each sheet (without changing which tab is active) ---
$wsSummary.Range("A4").Select()
$wsRepay.Range("G18:G19").Select()
$wsTrans.Range("D4").Select()

# Transactions was the active sheet before editing and must remain so
$wsTrans.Activate()
